# Weekly refresh of the Arveja Verde (Vega Monumental Concepcion) price series.
# Each data row (2-20) is updated in place with a newer weeks reported
# figures (date, variety, volume, min/max/weighted price, unit, origin, $/Kg).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44335    # D2 Fecha
$ws.Cells.Item(2, 8).Value = 'Perfection'    # H2 Variedad
$ws.Cells.Item(2, 10).Value = 100    # J2 Volumen
$ws.Cells.Item(2, 11).Value = 30000    # K2 Precio minimo
$ws.Cells.Item(2, 12).Value = 32000    # L2 Precio maximo
$ws.Cells.Item(2, 13).Value = 31000    # M2 Precio promedio ponderado
$ws.Cells.Item(2, 14).Value = '$/malla 25 kilos'    # N2 Unidad de comercializacion
$ws.Cells.Item(2, 15).Value = 'Provincia de Huasco'    # O2 Origen
$ws.Cells.Item(2, 16).Value = 1240    # P2 Precio $/Kg

$ws.Cells.Item(3, 4).Value = 44503    # D3 Fecha
$ws.Cells.Item(3, 8).Value = 'Perfection'    # H3 Variedad
$ws.Cells.Item(3, 10).Value = 200    # J3 Volumen
$ws.Cells.Item(3, 11).Value = 15000    # K3 Precio minimo
$ws.Cells.Item(3, 12).Value = 16000    # L3 Precio maximo
$ws.Cells.Item(3, 13).Value = 15500    # M3 Precio promedio ponderado
$ws.Cells.Item(3, 14).Value = '$/malla 25 kilos'    # N3 Unidad de comercializacion
$ws.Cells.Item(3, 15).Value = 'Provincia de Limarí'    # O3 Origen
$ws.Cells.Item(3, 16).Value = 620    # P3 Precio $/Kg

$ws.Cells.Item(4, 4).Value = 44399    # D4 Fecha
$ws.Cells.Item(4, 8).Value = 'Perfection'    # H4 Variedad
$ws.Cells.Item(4, 10).Value = 50    # J4 Volumen
$ws.Cells.Item(4, 11).Value = 39000    # K4 Precio minimo
$ws.Cells.Item(4, 12).Value = 40000    # L4 Precio maximo
$ws.Cells.Item(4, 13).Value = 39600    # M4 Precio promedio ponderado
$ws.Cells.Item(4, 14).Value = '$/malla 25 kilos'    # N4 Unidad de comercializacion
$ws.Cells.Item(4, 15).Value = 'Provincia de Huasco'    # O4 Origen
$ws.Cells.Item(4, 16).Value = 1584    # P4 Precio $/Kg

$ws.Cells.Item(5, 4).Value = 44496    # D5 Fecha
$ws.Cells.Item(5, 8).Value = 'Perfection'    # H5 Variedad
$ws.Cells.Item(5, 10).Value = 250    # J5 Volumen
$ws.Cells.Item(5, 11).Value = 14000    # K5 Precio minimo
$ws.Cells.Item(5, 12).Value = 15000    # L5 Precio maximo
$ws.Cells.Item(5, 13).Value = 14520    # M5 Precio promedio ponderado
$ws.Cells.Item(5, 14).Value = '$/malla 25 kilos'    # N5 Unidad de comercializacion
$ws.Cells.Item(5, 15).Value = 'Provincia de Huasco'    # O5 Origen
$ws.Cells.Item(5, 16).Value = 581    # P5 Precio $/Kg

$ws.Cells.Item(6, 4).Value = 44342    # D6 Fecha
$ws.Cells.Item(6, 8).Value = 'Perfection'    # H6 Variedad
$ws.Cells.Item(6, 10).Value = 60    # J6 Volumen
$ws.Cells.Item(6, 11).Value = 30000    # K6 Precio minimo
$ws.Cells.Item(6, 12).Value = 32000    # L6 Precio maximo
$ws.Cells.Item(6, 13).Value = 31000    # M6 Precio promedio ponderado
$ws.Cells.Item(6, 14).Value = '$/malla 25 kilos'    # N6 Unidad de comercializacion
$ws.Cells.Item(6, 15).Value = 'Provincia de Limarí'    # O6 Origen
$ws.Cells.Item(6, 16).Value = 1240    # P6 Precio $/Kg

$ws.Cells.Item(7, 4).Value = 44517    # D7 Fecha
$ws.Cells.Item(7, 8).Value = 'Perfection'    # H7 Variedad
$ws.Cells.Item(7, 10).Value = 110    # J7 Volumen
$ws.Cells.Item(7, 11).Value = 17000    # K7 Precio minimo
$ws.Cells.Item(7, 12).Value = 18000    # L7 Precio maximo
$ws.Cells.Item(7, 13).Value = 17455    # M7 Precio promedio ponderado
$ws.Cells.Item(7, 14).Value = '$/saco 25 kilos'    # N7 Unidad de comercializacion
$ws.Cells.Item(7, 15).Value = 'Región del Maule'    # O7 Origen
$ws.Cells.Item(7, 16).Value = 698    # P7 Precio $/Kg

$ws.Cells.Item(8, 4).Value = 44512    # D8 Fecha
$ws.Cells.Item(8, 8).Value = 'Perfection'    # H8 Variedad
$ws.Cells.Item(8, 10).Value = 100    # J8 Volumen
$ws.Cells.Item(8, 11).Value = 14000    # K8 Precio minimo
$ws.Cells.Item(8, 12).Value = 15000    # L8 Precio maximo
$ws.Cells.Item(8, 13).Value = 14500    # M8 Precio promedio ponderado
$ws.Cells.Item(8, 14).Value = '$/saco 25 kilos'    # N8 Unidad de comercializacion
$ws.Cells.Item(8, 15).Value = 'Región del Maule'    # O8 Origen
$ws.Cells.Item(8, 16).Value = 580    # P8 Precio $/Kg

$ws.Cells.Item(9, 4).Value = 44483    # D9 Fecha
$ws.Cells.Item(9, 8).Value = 'Perfection'    # H9 Variedad
$ws.Cells.Item(9, 10).Value = 220    # J9 Volumen
$ws.Cells.Item(9, 11).Value = 19000    # K9 Precio minimo
$ws.Cells.Item(9, 12).Value = 20000    # L9 Precio maximo
$ws.Cells.Item(9, 13).Value = 19455    # M9 Precio promedio ponderado
$ws.Cells.Item(9, 14).Value = '$/saco 25 kilos'    # N9 Unidad de comercializacion
$ws.Cells.Item(9, 15).Value = 'Región Metropolitana'    # O9 Origen
$ws.Cells.Item(9, 16).Value = 778    # P9 Precio $/Kg

$ws.Cells.Item(10, 4).Value = 44505    # D10 Fecha
$ws.Cells.Item(10, 8).Value = 'Perfection'    # H10 Variedad
$ws.Cells.Item(10, 10).Value = 210    # J10 Volumen
$ws.Cells.Item(10, 11).Value = 6500    # K10 Precio minimo
$ws.Cells.Item(10, 12).Value = 7000    # L10 Precio maximo
$ws.Cells.Item(10, 13).Value = 6714    # M10 Precio promedio ponderado
$ws.Cells.Item(10, 14).Value = '$/malla 25 kilos'    # N10 Unidad de comercializacion
$ws.Cells.Item(10, 15).Value = 'Región del Maule'    # O10 Origen
$ws.Cells.Item(10, 16).Value = 269    # P10 Precio $/Kg

$ws.Cells.Item(11, 4).Value = 44162    # D11 Fecha
$ws.Cells.Item(11, 8).Value = 'Sin especificar'    # H11 Variedad
$ws.Cells.Item(11, 10).Value = 100    # J11 Volumen
$ws.Cells.Item(11, 11).Value = 17000    # K11 Precio minimo
$ws.Cells.Item(11, 12).Value = 18000    # L11 Precio maximo
$ws.Cells.Item(11, 13).Value = 17500    # M11 Precio promedio ponderado
$ws.Cells.Item(11, 14).Value = '$/saco 25 kilos'    # N11 Unidad de comercializacion
$ws.Cells.Item(11, 15).Value = 'Región del Maule'    # O11 Origen
$ws.Cells.Item(11, 16).Value = 700    # P11 Precio $/Kg

$ws.Cells.Item(12, 4).Value = 44482    # D12 Fecha
$ws.Cells.Item(12, 8).Value = 'Perfection'    # H12 Variedad
$ws.Cells.Item(12, 10).Value = 130    # J12 Volumen
$ws.Cells.Item(12, 11).Value = 24000    # K12 Precio minimo
$ws.Cells.Item(12, 12).Value = 25000    # L12 Precio maximo
$ws.Cells.Item(12, 13).Value = 24385    # M12 Precio promedio ponderado
$ws.Cells.Item(12, 14).Value = '$/saco 25 kilos'    # N12 Unidad de comercializacion
$ws.Cells.Item(12, 15).Value = 'Región de O''Higgins'    # O12 Origen
$ws.Cells.Item(12, 16).Value = 975    # P12 Precio $/Kg

$ws.Cells.Item(13, 4).Value = 44533    # D13 Fecha
$ws.Cells.Item(13, 8).Value = 'Perfection'    # H13 Variedad
$ws.Cells.Item(13, 10).Value = 80    # J13 Volumen
$ws.Cells.Item(13, 11).Value = 14000    # K13 Precio minimo
$ws.Cells.Item(13, 12).Value = 15000    # L13 Precio maximo
$ws.Cells.Item(13, 13).Value = 14375    # M13 Precio promedio ponderado
$ws.Cells.Item(13, 14).Value = '$/malla 25 kilos'    # N13 Unidad de comercializacion
$ws.Cells.Item(13, 15).Value = 'Región del Maule'    # O13 Origen
$ws.Cells.Item(13, 16).Value = 575    # P13 Precio $/Kg

$ws.Cells.Item(14, 4).Value = 44539    # D14 Fecha
$ws.Cells.Item(14, 8).Value = 'Sin especificar'    # H14 Variedad
$ws.Cells.Item(14, 10).Value = 50    # J14 Volumen
$ws.Cells.Item(14, 11).Value = 13000    # K14 Precio minimo
$ws.Cells.Item(14, 12).Value = 14000    # L14 Precio maximo
$ws.Cells.Item(14, 13).Value = 13400    # M14 Precio promedio ponderado
$ws.Cells.Item(14, 14).Value = '$/saco 25 kilos'    # N14 Unidad de comercializacion
$ws.Cells.Item(14, 15).Value = 'Región del Maule'    # O14 Origen
$ws.Cells.Item(14, 16).Value = 536    # P14 Precio $/Kg

$ws.Cells.Item(15, 4).Value = 44532    # D15 Fecha
$ws.Cells.Item(15, 8).Value = 'Sin especificar'    # H15 Variedad
$ws.Cells.Item(15, 10).Value = 250    # J15 Volumen
$ws.Cells.Item(15, 11).Value = 14000    # K15 Precio minimo
$ws.Cells.Item(15, 12).Value = 15000    # L15 Precio maximo
$ws.Cells.Item(15, 13).Value = 14400    # M15 Precio promedio ponderado
$ws.Cells.Item(15, 14).Value = '$/saco 25 kilos'    # N15 Unidad de comercializacion
$ws.Cells.Item(15, 15).Value = 'Región del Maule'    # O15 Origen
$ws.Cells.Item(15, 16).Value = 576    # P15 Precio $/Kg

$ws.Cells.Item(16, 4).Value = 44328    # D16 Fecha
$ws.Cells.Item(16, 8).Value = 'Perfection'    # H16 Variedad
$ws.Cells.Item(16, 10).Value = 100    # J16 Volumen
$ws.Cells.Item(16, 11).Value = 33000    # K16 Precio minimo
$ws.Cells.Item(16, 12).Value = 34000    # L16 Precio maximo
$ws.Cells.Item(16, 13).Value = 33500    # M16 Precio promedio ponderado
$ws.Cells.Item(16, 14).Value = '$/malla 25 kilos'    # N16 Unidad de comercializacion
$ws.Cells.Item(16, 15).Value = 'Provincia de Huasco'    # O16 Origen
$ws.Cells.Item(16, 16).Value = 1340    # P16 Precio $/Kg

$ws.Cells.Item(17, 4).Value = 44519    # D17 Fecha
$ws.Cells.Item(17, 8).Value = 'Perfection'    # H17 Variedad
$ws.Cells.Item(17, 10).Value = 240    # J17 Volumen
$ws.Cells.Item(17, 11).Value = 17000    # K17 Precio minimo
$ws.Cells.Item(17, 12).Value = 18000    # L17 Precio maximo
$ws.Cells.Item(17, 13).Value = 17583    # M17 Precio promedio ponderado
$ws.Cells.Item(17, 14).Value = '$/saco 25 kilos'    # N17 Unidad de comercializacion
$ws.Cells.Item(17, 15).Value = 'Carahue'    # O17 Origen
$ws.Cells.Item(17, 16).Value = 703    # P17 Precio $/Kg

$ws.Cells.Item(18, 4).Value = 44454    # D18 Fecha
$ws.Cells.Item(18, 8).Value = 'Perfection'    # H18 Variedad
$ws.Cells.Item(18, 10).Value = 100    # J18 Volumen
$ws.Cells.Item(18, 11).Value = 36000    # K18 Precio minimo
$ws.Cells.Item(18, 12).Value = 38000    # L18 Precio maximo
$ws.Cells.Item(18, 13).Value = 37000    # M18 Precio promedio ponderado
$ws.Cells.Item(18, 14).Value = '$/malla 25 kilos'    # N18 Unidad de comercializacion
$ws.Cells.Item(18, 15).Value = 'Provincia de Limarí'    # O18 Origen
$ws.Cells.Item(18, 16).Value = 1480    # P18 Precio $/Kg

$ws.Cells.Item(19, 4).Value = 44518    # D19 Fecha
$ws.Cells.Item(19, 8).Value = 'Perfection'    # H19 Variedad
$ws.Cells.Item(19, 10).Value = 350    # J19 Volumen
$ws.Cells.Item(19, 11).Value = 14000    # K19 Precio minimo
$ws.Cells.Item(19, 12).Value = 15000    # L19 Precio maximo
$ws.Cells.Item(19, 13).Value = 14571    # M19 Precio promedio ponderado
$ws.Cells.Item(19, 14).Value = '$/saco 25 kilos'    # N19 Unidad de comercializacion
$ws.Cells.Item(19, 15).Value = 'Región del Maule'    # O19 Origen
$ws.Cells.Item(19, 16).Value = 583    # P19 Precio $/Kg

$ws.Cells.Item(20, 4).Value = 44540    # D20 Fecha
$ws.Cells.Item(20, 8).Value = 'Sin especificar'    # H20 Variedad
$ws.Cells.Item(20, 10).Value = 110    # J20 Volumen
$ws.Cells.Item(20, 11).Value = 16000    # K20 Precio minimo
$ws.Cells.Item(20, 12).Value = 17000    # L20 Precio maximo
$ws.Cells.Item(20, 13).Value = 31000    # M20 Precio promedio ponderado
$ws.Cells.Item(20, 14).Value = '$/saco 25 kilos'    # N20 Unidad de comercializacion
$ws.Cells.Item(20, 15).Value = 'Región del Maule'    # O20 Origen
$ws.Cells.Item(20, 16).Value = 662    # P20 Precio $/Kg

